$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2332155477031802
$ws.Range("C2").Value = 0.4840989399293286
$ws.Range("P2").Value = 0.1872791519434629
$ws.Range("S2").Value = 0.09540636042402827
$ws.Range("B3").Value = 0.01470588235294118
$ws.Range("C3").Value = 0.02941176470588235
$ws.Range("J3").Value = 0.007352941176470588
$ws.Range("O3").Value = 0.007352941176470588
$ws.Range("P3").Value = 0.7132352941176471
$ws.Range("S3").Value = 0.2279411764705882
$ws.Range("J4").Value = 0.0425531914893617
$ws.Range("P4").Value = 0.6170212765957447
$ws.Range("S4").Value = 0.3404255319148936
$ws.Range("B6").Value = 0.06806282722513089
$ws.Range("D6").Value = 0.03664921465968586
$ws.Range("F6").Value = 0.0418848167539267
$ws.Range("J6").Value = 0.256544502617801
$ws.Range("O6").Value = 0.01570680628272251
$ws.Range("Q6").Value = 0.1780104712041885
$ws.Range("R6").Value = 0.0418848167539267
$ws.Range("S6").Value = 0.3612565445026178
$ws.Range("B7").Value = 0.08187134502923976
$ws.Range("D7").Value = 0.04678362573099415
$ws.Range("F7").Value = 0.05263157894736842
$ws.Range("J7").Value = 0.1871345029239766
$ws.Range("O7").Value = 0.005847953216374269
$ws.Range("Q7").Value = 0.1111111111111111
$ws.Range("R7").Value = 0.06432748538011696
$ws.Range("S7").Value = 0.4502923976608187
$ws.Range("B8").Value = 0.087890625
$ws.Range("D8").Value = 0.015625
$ws.Range("F8").Value = 0.068359375
$ws.Range("J8").Value = 0.125
$ws.Range("O8").Value = 0.017578125
$ws.Range("Q8").Value = 0.158203125
$ws.Range("R8").Value = 0.095703125
$ws.Range("S8").Value = 0.431640625
$ws.Range("B9").Value = 0.06951871657754011
$ws.Range("D9").Value = 0.0106951871657754
$ws.Range("E9").Value = 0.0053475935828877
$ws.Range("F9").Value = 0.053475935828877
$ws.Range("J9").Value = 0.1550802139037433
$ws.Range("O9").Value = 0.0106951871657754
$ws.Range("Q9").Value = 0.1390374331550802
$ws.Range("R9").Value = 0.09625668449197861
$ws.Range("S9").Value = 0.4598930481283423
$ws.Range("B10").Value = 0.1147540983606557
$ws.Range("D10").Value = 0.02157031924072476
$ws.Range("F10").Value = 0.06212251941328732
$ws.Range("J10").Value = 0.1363244176013805
$ws.Range("O10").Value = 0.01035375323554789
$ws.Range("Q10").Value = 0.1949956859361519
$ws.Range("R10").Value = 0.06902502157031924
$ws.Range("S10").Value = 0.3908541846419327
$ws.Range("G11").Value = 0.15
$ws.Range("J11").Value = 0.1071428571428571
$ws.Range("K11").Value = 0.2321428571428572
$ws.Range("L11").Value = 0.4928571428571429
$ws.Range("S11").Value = 0.01785714285714286
$ws.Range("G12").Value = 0.773972602739726
$ws.Range("J12").Value = 0.1438356164383562
$ws.Range("K12").Value = 0.00684931506849315
$ws.Range("L12").Value = 0.0410958904109589
$ws.Range("S12").Value = 0.03424657534246575
$ws.Range("G13").Value = 0.5833333333333334
$ws.Range("J13").Value = 0.3611111111111111
$ws.Range("S13").Value = 0.05555555555555555
$ws.Range("F14").Value = 0.5
$ws.Range("J14").Value = 0.5
$ws.Range("F15").Value = 0.01744186046511628
$ws.Range("H15").Value = 0.1976744186046512
$ws.Range("I15").Value = 0.0872093023255814
$ws.Range("J15").Value = 0.3197674418604651
$ws.Range("K15").Value = 0.05232558139534884
$ws.Range("M15").Value = 0.005813953488372093
$ws.Range("O15").Value = 0.01744186046511628
$ws.Range("S15").Value = 0.3023255813953488
$ws.Range("F16").Value = 0.02298850574712644
$ws.Range("H16").Value = 0.1781609195402299
$ws.Range("I16").Value = 0.1379310344827586
$ws.Range("J16").Value = 0.3505747126436782
$ws.Range("K16").Value = 0.1091954022988506
$ws.Range("O16").Value = 0.04597701149425287
$ws.Range("S16").Value = 0.1551724137931035
$ws.Range("F17").Value = 0.02077922077922078
$ws.Range("H17").Value = 0.1922077922077922
$ws.Range("I17").Value = 0.09350649350649351
$ws.Range("J17").Value = 0.3792207792207792
$ws.Range("K17").Value = 0.09870129870129871
$ws.Range("M17").Value = 0.01298701298701299
$ws.Range("N17").Value = 0.002597402597402597
$ws.Range("O17").Value = 0.04415584415584416
$ws.Range("S17").Value = 0.1558441558441558
$ws.Range("F18").Value = 0.01204819277108434
$ws.Range("H18").Value = 0.2048192771084337
$ws.Range("I18").Value = 0.0963855421686747
$ws.Range("J18").Value = 0.4096385542168675
$ws.Range("K18").Value = 0.06626506024096386
$ws.Range("M18").Value = 0.01204819277108434
$ws.Range("O18").Value = 0.006024096385542169
$ws.Range("S18").Value = 0.1927710843373494
$ws.Range("F19").Value = 0.01255539143279173
$ws.Range("H19").Value = 0.2518463810930576
$ws.Range("I19").Value = 0.07385524372230429
$ws.Range("J19").Value = 0.3293943870014771
$ws.Range("K19").Value = 0.0982274741506647
$ws.Range("M19").Value = 0.01994091580502216
$ws.Range("N19").Value = 0.0007385524372230429
$ws.Range("O19").Value = 0.06573116691285082
$ws.Range("S19").Value = 0.1477104874446086

Write-Host "Applied 110 changes"
